$wb = $excel.ActiveWorkbook

# Rename the original sheet to "road"
$road = $wb.Worksheets.Item(1)
$road.Name = "road"

# Add a new sheet "building" right after "road" and copy the data into it
$building = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $road)
$building.Name = "building"

[void]$road.Range("A1:C9").Copy()
[void]$building.Range("A1").PasteSpecial()

# Restore the per-sheet selections
[void]$road.Range("A1:C12").Select()

[void]$building.Range("E17").Select()
[void]$building.Activate()
